$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.553.61"
$ws.Range("E2").Value = "'  -0.18%  "
$ws.Range("D3").Value = "'1.729.39"
$ws.Range("E3").Value = "'  -0.84%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'246.07"
$ws.Range("E5").Value = "'  -0.52%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  +0.00%  "
$ws.Range("E7").Value = "'  +0.42%  "
$ws.Range("E8").Value = "'  -1.04%  "
$ws.Range("D9").Value = "'0.06226"
$ws.Range("E9").Value = "'  -0.34%  "
$ws.Range("D10").Value = "'1.733.54"
$ws.Range("E10").Value = "'  -0.58%  "
$ws.Range("D11").Value = "'0.07094"
$ws.Range("E11").Value = "'  -0.42%  "
$ws.Range("D12").Value = "'15.57"
$ws.Range("E12").Value = "'  -1.43%  "
$ws.Range("D13").Value = "'0.6082"
$ws.Range("E13").Value = "'  -1.80%  "
$ws.Range("D14").Value = "'4.549"
$ws.Range("E14").Value = "'  +1.01%  "
$ws.Range("D15").Value = "'77.23"
$ws.Range("E15").Value = "'  -0.47%  "
$ws.Range("E16").Value = "'  +0.04%  "
$ws.Range("D17").Value = "'26.545.03"
$ws.Range("E17").Value = "'  -0.21%  "
$ws.Range("B18").Value = "'ShibaInu"
$ws.Range("C18").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007295"
$ws.Range("E18").Value = "'  +5.79%  "
$ws.Range("B19").Value = "'BinanceUSD"
$ws.Range("C19").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "'  +0.00%  "
$ws.Range("E20").Value = "'  -1.92%  "
$ws.Range("D21").Value = "'1.954.51"
$ws.Range("E21").Value = "'  -0.54%  "
$ws.Range("E22").Value = "'  -2.89%  "
$ws.Range("D23").Value = "'8.769"
$ws.Range("E23").Value = "'  -0.72%  "
$ws.Range("D24").Value = "'5.238"
$ws.Range("E24").Value = "'  -2.08%  "
$ws.Range("D25").Value = "'137.33"
$ws.Range("D26").Value = "'15.44"
$ws.Range("E26").Value = "'  +0.18%  "
$ws.Range("B27").Value = "'LidoDAOToken"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'1.774"
$ws.Range("E27").Value = "'  -2.38%  "
$ws.Range("B28").Value = "'Toncoin"
$ws.Range("C28").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'1.406"
$ws.Range("E28").Value = "'  -1.94%  "
$ws.Range("D29").Value = "'108.18"
$ws.Range("E29").Value = "'  +0.84%  "
$ws.Range("D30").Value = "'3.959"
$ws.Range("E30").Value = "'  -1.40%  "
$ws.Range("D31").Value = "'0.08004"
$ws.Range("E31").Value = "'  +1.48%  "
$ws.Range("D32").Value = "'3.692"
$ws.Range("E32").Value = "'  -1.69%  "
$ws.Range("D33").Value = "'0.04566"
$ws.Range("E33").Value = "'  -0.44%  "
$ws.Range("D34").Value = "'0.9999"
$ws.Range("E34").Value = "'  -0.01%  "
$ws.Range("E35").Value = "'  +0.11%  "
$ws.Range("D36").Value = "'0.9998"
$ws.Range("E36").Value = "'  +0.15%  "
$ws.Range("D37").Value = "'0.6306"
$ws.Range("E37").Value = "'  -1.43%  "
$ws.Range("D38").Value = "'0.8932"
$ws.Range("E38").Value = "'  -5.32%  "
$ws.Range("D39").Value = "'2.001"
$ws.Range("E39").Value = "'  +0.23%  "
$ws.Range("E40").Value = "'  -1.86%  "
$ws.Range("E41").Value = "'  +0.13%  "
$ws.Range("D42").Value = "'0.01502"
$ws.Range("E42").Value = "'  -0.35%  "
$ws.Range("D43").Value = "'101.69"
$ws.Range("E43").Value = "'  -10.24%  "
$ws.Range("D44").Value = "'5.386"
$ws.Range("E44").Value = "'  -6.39%  "
$ws.Range("D45").Value = "'0.3886"
$ws.Range("E45").Value = "'  -0.92%  "
$ws.Range("D46").Value = "'7.054"
$ws.Range("E46").Value = "'  +5.47%  "
$ws.Range("D47").Value = "'0.1181"
$ws.Range("E47").Value = "'  -2.12%  "
$ws.Range("D48").Value = "'0.05393"
$ws.Range("E48").Value = "'  +1.18%  "
$ws.Range("D49").Value = "'7.918"
$ws.Range("E49").Value = "'  -0.11%  "
$ws.Range("D50").Value = "'30.58"
$ws.Range("E50").Value = "'  -0.82%  "
$ws.Range("D51").Value = "'1.256"
$ws.Range("E51").Value = "'  -0.84%  "
